$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview sheet.
# This text is shared with de-de!H2 ("Correspond Handoff Datetime") which
# carried the same original timestamp, so both must be updated together.
$wsOverview.Range("G2").Value = "2016-09-01 01:09:13"
$wsDeDe.Range("H2").Value = "2016-09-01 01:09:13"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-01 01:09:07"
$wsZhCn.Range("K2").Value = "2016-09-01 01:09:31"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-01 01:09:40"
